$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EREC row (row 8): stock balance text "6:8" -> "7:9"
$ws.Range("H8").Value = "7:9"

# NETLOOK row (row 9): stock balance text "1:0" -> "2:0"
$ws.Range("H9").Value = "2:0"

# Footer timestamp: "Wednesday, 30 July, 2025 12:32 AM" -> "...12:33 AM"
$ws.Range("A11").Value = "Wednesday, 30 July, 2025 12:33 AM"
